{"js": "// Apply 5 text replacements inside table cells describing the\n// \"Registro de usuario\" use case, per the commit's final edits.\n\nconst replacements = [\n  {\n    find: \"El sistema desplegar\u00e1 en m\u00f3vil la geolocalizaci\u00f3n y en web el cat\u00e1logo de productos.\",\n    replace: \"El sistema desplegar\u00e1 el men\u00fa de inicio\"\n  },\n  {\n    find: \"El Proveedor selecciona el icono de perfil.\",\n    replace: \"El Proveedor selecciona perfil.\"\n  },\n  {\n    find: \"El sistema despliega la interfaz de \u00bfquieres vender tus productos?\",\n    replace: \"El sistema despliega la interfaz de registro\"\n  },\n  {\n    find: \"El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su correo electr\u00f3nico y contrase\u00f1a.\",\n    replace: \"El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su rut y contrase\u00f1a.\"\n  },\n  {\n    find: \"El correo electr\u00f3nico ingresado ya est\u00e1 registrado: el sistema muestra un mensaje de error indicando que el correo ya est\u00e1 en uso.\",\n    replace: \"El correo electr\u00f3nico o rut ingresado ya est\u00e1 registrado: el sistema muestra un mensaje de error indicando que el correo ya est\u00e1 en uso.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply 5 text replacements inside table cells describing the\n# \"Registro de usuario\" use case, per the commit's final edits.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"El sistema desplegar\u00e1 en m\u00f3vil la geolocalizaci\u00f3n y en web el cat\u00e1logo de productos.\"; Replace = \"El sistema desplegar\u00e1 el men\u00fa de inicio\" },\n    @{ Find = \"El Proveedor selecciona el icono de perfil.\"; Replace = \"El Proveedor selecciona perfil.\" },\n    @{ Find = \"El sistema despliega la interfaz de \u00bfquieres vender tus productos?\"; Replace = \"El sistema despliega la interfaz de registro\" },\n    @{ Find = \"El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su correo electr\u00f3nico y contrase\u00f1a.\"; Replace = \"El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su rut y contrase\u00f1a.\" },\n    @{ Find = \"El correo electr\u00f3nico ingresado ya est\u00e1 registrado: el sistema muestra un mensaje de error indicando que el correo ya est\u00e1 en uso.\"; Replace = \"El correo electr\u00f3nico o rut ingresado ya est\u00e1 registrado: el sistema muestra un mensaje de error indicando que el correo ya est\u00e1 en uso.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
